# "Stat server partially in place"
# Multiple login changes and corrections.
#
# - Insert a new row above the "game1"/"game2" port block and start
#   filling in the (still partial) "Stat" server entry, with the base
#   port struck through since it isn't live yet.
# - This pushes the existing game1 (row 11 -> 12) and game2 (row 12 -> 13)
#   rows down by one.
# - Leave the selection sitting on C11 (the still-blank cell below the
#   new Stat row) and flip the sheet to portrait for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at 10; game1/game2 shift down to 12/13.
$ws.Rows("10:10").Insert()

# New "Stat" server row. The base port (B10) is marked struck-through
# to show it's only partially configured so far.
$ws.Range("A10").Value = "Stat"
$ws.Range("B10").Value = 7800
$ws.Range("C10").Value = 7802
$ws.Range("B10").Font.Strikethrough = $true

# Sheet is set up for printing in portrait orientation.
$ws.PageSetup.Orientation = 1

# Leave the cursor where the edit left it.
$ws.Range("C11").Select()
